# Revert "adding term 2.0 now utf-8"
$wb = $excel.ActiveWorkbook

# --- Update values on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Update value on the "Include from FSIII" sheet ---
$inc1 = $wb.Worksheets.Item("Include from FSIII")
$inc1.Range("C2").Value = "F"

# --- Remove the duplicate "Include from FSIII 2" sheet ---
$inc2 = $wb.Worksheets.Item("Include from FSIII 2")
$inc2.Delete()
